$d = $word.ActiveDocument

# The paragraph currently reads:
#   "When the zip file is expanded, these are the projects that are to be run:"
# A new bold sentence is inserted at the very start of the paragraph:
#   "Expand the modt recent zip file. "
# and the hidden "_GoBack" bookmark - previously sitting right before the
# word "projects" - moves so it sits right after the newly inserted
# sentence (i.e. right before "When the zip file").

# Locate the paragraph that currently begins with "When the zip file" —
# this is the paragraph the edit targets.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith("When the zip file")) {
        $target = $candidate
        break
    }
}

$r = $target.Range
$start = $r.Start
$text = "Expand the modt recent zip file. "

# Insert the new bold sentence at the very start of the paragraph.
$ins = $d.Range($start, $start)
$ins.InsertBefore($text)
$newRange = $d.Range($start, $start + $text.Length)
$newRange.Bold = $true

# Relocate the "_GoBack" bookmark so it sits right after the inserted
# sentence (it previously sat right before "projects").
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}
$bmPos = $start + $text.Length
$bmTarget = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmTarget)

Write-Output $target.Range.Text
